# WIP: translating rolling window to python
# Update the computed metrics for each portfolio (rows 2-4, cols B-F)
# with refreshed values from the rolling-window calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Equal Weights
$ws.Range("B2").Value = 17.66844471405463
$ws.Range("C2").Value = 23.73070317375016
$ws.Range("D2").Value = 0.7420110809667156
$ws.Range("E2").Value = 1.305573997661102
$ws.Range("F2").Value = 6.992209217887703

# Row 3: HRP
$ws.Range("B3").Value = 14.73541984188003
$ws.Range("C3").Value = 18.65086695561742
$ws.Range("D3").Value = 0.7868492052837236
$ws.Range("E3").Value = 1.209939020860741
$ws.Range("F3").Value = 5.495682552702294

# Row 4: MV
$ws.Range("B4").Value = 13.05291861646707
$ws.Range("C4").Value = 17.09425770168847
$ws.Range("D4").Value = 0.760075040590017
$ws.Range("E4").Value = 1.169449336231703
$ws.Range("F4").Value = 5.022278589382549
